$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder "Recorded By" list (G2) ---
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 3: reorder "Recorded By" list (G3) ---
$ws.Range("G3").Value = "System, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 4: reorder "Recorded By" list (G4) ---
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 5: reorder "Recorded By" list (G5) ---
$ws.Range("G5").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 6: reorder "Recorded By" list (G6) ---
$ws.Range("G6").Value = "alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# --- Row 7: reorder "Recorded By" list (G7) + Missing Sessions stat (L7) ---
$ws.Range("G7").Value = "NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("L7").Value = 3

# --- Row 8: Pending Sessions stat (L8) ---
$ws.Range("L8").Value = 1

# --- Row 12: reorder "Recorded By" list (G12) ---
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

# --- Row 13: recolor to match "Not Recorded" rows (copy format from row 11) and update status text ---
$ws.Range("A11:I11").Copy()
$ws.Range("A13:I13").PasteSpecial(-4122)
$ws.Range("I13").Value = "Not Recorded"

# --- Row 15: reorder "Recorded By" list (G15) + Missing/Pending stats (P15/Q15) ---
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 1

# --- Row 17: reorder "Recorded By" list (G17) ---
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 20: reorder "Recorded By" list (G20) ---
$ws.Range("G20").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 24: reorder "Recorded By" list (G24) ---
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"

# --- Row 25: reorder "Recorded By" list (G25) ---
$ws.Range("G25").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"

# --- Row 27: reorder "Recorded By" list (G27) ---
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 30: reorder "Recorded By" list (G30) ---
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
